$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.445.96"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "2.294.33"
$ws.Range("E3").Value = "  -0.03%  "
$cell = $ws.Range("D4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = $origStyle
$ws.Range("E4").Value = "  +0.02%  "
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "300.52"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  -2.21%  "
$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "94.67"
$cell.Style = $origStyle
$ws.Range("E6").Value = "  -1.81%  "
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").Value = "  +0.06%  "
$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.490"
$cell.Style = $origStyle
$ws.Range("E9").Value = "  -2.27%  "
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "34.33"
$cell.Style = $origStyle
$ws.Range("E10").Value = "  -2.63%  "
$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "18.96"
$cell.Style = $origStyle
$ws.Range("E11").Value = "  +3.05%  "
$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0778"
$cell.Style = $origStyle
$ws.Range("E12").Value = "  -1.47%  "
$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.69"
$cell.Style = $origStyle
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("D15").Value = "2.649.95"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "2.278.56"
$ws.Range("E16").Value = "  -0.73%  "
$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.777"
$cell.Style = $origStyle
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").Value = "42.402.35"
$ws.Range("E18").Value = "  -0.75%  "
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "12.15"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  -6.07%  "
$ws.Range("D20").Value = "0.0₃0886"
$ws.Range("E20").Value = "  -1.09%  "
$ws.Range("E21").Value = "  -1.52%  "
$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "67.47"
$cell.Style = $origStyle
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "235.40"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.25"
$cell.Style = $origStyle
$ws.Range("E24").Value = "  +5.38%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  -2.41%  "
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "24.19"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  -3.52%  "
$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.05"
$cell.Style = $origStyle
$ws.Range("E28").Value = "  -13.90%  "
$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "164.32"
$cell.Style = $origStyle
$ws.Range("E29").Value = "  -1.21%  "
$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.01"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  -0.42%  "
$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "31.53"
$cell.Style = $origStyle
$ws.Range("E31").Value = "  -4.26%  "
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = $origStyle
$ws.Range("E32").Value = "  +0.00%  "
$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.96"
$cell.Style = $origStyle
$ws.Range("E33").Value = "  -0.39%  "
$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "17.42"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  -1.01%  "
$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0693"
$cell.Style = $origStyle
$ws.Range("E35").Value = "  -0.02%  "
$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.32"
$cell.Style = $origStyle
$ws.Range("E36").Value = "  -3.11%  "
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.33"
$cell.Style = $origStyle
$ws.Range("E37").Value = "  -8.70%  "
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0997"
$cell.Style = $origStyle
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("E40").Value = "  -1.45%  "
$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.67"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  -1.12%  "
$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "19.81"
$cell.Style = $origStyle
$ws.Range("E42").Value = "  +8.78%  "
$ws.Range("D43").Value = "1.946.60"
$ws.Range("E43").Value = "  -3.12%  "
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "10.32"
$cell.Style = $origStyle
$ws.Range("E44").Value = "  +2.38%  "
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("E46").Value = "  +1.51%  "
$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.72"
$cell.Style = $origStyle
$ws.Range("E47").Value = "  -2.88%  "
$ws.Range("D48").Value = "2.519.24"
$ws.Range("E48").Value = "  -0.12%  "
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "52.70"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  -2.31%  "
$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.79"
$cell.Style = $origStyle
$ws.Range("E50").Value = "  -3.90%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.12"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  +0.41%  "
